$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 204, shifting existing rows 204:249 down to 205:250
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row 204 with the new price-report entry
$ws.Range("A204").Value = 3
$ws.Range("B204").Value = "Femacal de La Calera"
$ws.Range("C204").Value = "Coquimbo"
$ws.Range("D204").Value = 44511
$ws.Range("E204").Value = 5
$ws.Range("F204").Value = 100112003
$ws.Range("G204").Value = "Ajo"
$ws.Range("H204").Value = "Chino"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 78
$ws.Range("K204").Value = 16000
$ws.Range("L204").Value = 16500
$ws.Range("M204").Value = 16244
$ws.Range("N204").Value = "$/caja 10 kilos"
$ws.Range("O204").Value = "China"
$ws.Range("P204").Value = 1624
$ws.Range("Q204").Value = 10
$ws.Range("R204").Value = "Hortaliza"
